# Excel erweitert, import angepasst
# Adds two new number columns (nummer4, nummer5) and two new data rows
# to the existing "Nummerneinteilung" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers for the two new columns -------------------------------------
$ws.Range("D1").Value = "nummer4"
$ws.Range("E1").Value = "nummer5"

# --- New data for existing rows (columns D & E) ---------------------------
$ws.Range("D2").Value = 20220305
$ws.Range("E2").Value = 20220302

$ws.Range("D3").Value = 20220332
$ws.Range("E3").Value = 20220312

$ws.Range("D4").Value = 20220324
$ws.Range("E4").Value = 20220326

# --- Two additional rows across all five columns ---------------------------
$ws.Range("A5").Value = 20220405
$ws.Range("B5").Value = 20220408
$ws.Range("C5").Value = 20220405
$ws.Range("D5").Value = 20220405
$ws.Range("E5").Value = 20220401

$ws.Range("A6").Value = 20220501
$ws.Range("B6").Value = 20220507
$ws.Range("C6").Value = 20220509
$ws.Range("D6").Value = 20220501
$ws.Range("E6").Value = 20220501

# --- Match the column width of column C (19.42578125) on the new columns --
$ws.Columns("D:E").ColumnWidth = $ws.Columns("C").ColumnWidth

# --- Move the selection to C7, matching where the user ended up next ------
$ws.Activate() | Out-Null
$ws.Range("C7").Select() | Out-Null

# --- Reposition the workbook window (second monitor / wider layout) -------
$win = $wb.Windows.Item(1)
$win.Left = -28920
$win.Top = -120
